$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ticket Sheet")

# Column O was recalculated off ticket quantity (I) instead of the total (L).
# O2 holds its own (non-shared) formula; O3:O10 are a shared formula group.
$ws.Range("O2").Formula = "=I2*0.2"
$ws.Range("O3:O10").Formula = "=I3*0.2"

# The extra catering/ticket option recorded in H6 was removed.
$ws.Range("H6").ClearContents()

# A new blank row was added below the totals, formatted like the currency cells.
$ws.Range("D15").Value = $null
$ws.Range("D15").NumberFormat = $ws.Range("O2").NumberFormat

# View state: zoom in to 140% and move the active selection to J14.
$excel.ActiveWindow.Zoom = 140
$ws.Range("J14").Select() | Out-Null
